# Apply the edits described by the diff:
# - D8 text: move the "Enfin il est courant..." sentence from the middle to the end of the paragraph
# - Sheet view: zoom 143% -> 70%, selection D8 -> D9
# - Row heights: minor adjustments on rows 6, 8, 9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E3-1-ISI1-4-005 et 6 SCEN")
$ws.Activate()

$d8Text = @"
Il est courant de représenter le SI par un modèle en 5 couches selon le principe de séparation des préocupations.
Chaque couche isole un aspect particulier du système d'information en étant responsable des intéractions entre ses éléments.
Chaque couche n'échangeant qu'avec ses couches adjacentes.
Les deux premières couches forment le système informatique, l'ensemble structuré des composants matériels et logiciels et les données permettant d'automatiser tout ou partie du système métier au travers de fonctionnalités qui lui sont nécessaires.
Le système métier est formé des services et processus de l'entreprise, des organisations qui les mettent en œuvre et des objets métier associés.
Un objet métier est un concept ou une abstraction ayant un sens pour des acteurs (partie prenante interne) d'une organisation (par exemple une entreprise). L'objet métier permet de décrire les entités manipulées par les acteurs dans le cadre de la description du métier.
Enfin il est courant d'y ajouter la vue stratégie. Il s'agit de la stratégie décidée par la D.S.I. pour le S.I. en alignement avec la stratégie de l'organisation.
"@
$ws.Range("D8").Value = $d8Text

# Sheet view: zoom + active selection
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("D9").Select()

# Row height tweaks
$ws.Rows.Item(6).RowHeight = 154.05
$ws.Rows.Item(8).RowHeight = 253.95
$ws.Rows.Item(9).RowHeight = 235.95
